$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; this shifts the existing rows 52-157 down
# to 53-158 (and Excel extends the used range / dimension accordingly).
$ws.Rows(52).Insert()

# Populate the newly inserted row 52 with the new weekly record.
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value = 44953
$ws.Cells.Item(52, 5).Value = 15
$ws.Cells.Item(52, 6).Value = 100112042
$ws.Cells.Item(52, 7).Value = "Locoto"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 190
$ws.Cells.Item(52, 11).Value = 19000
$ws.Cells.Item(52, 12).Value = 20000
$ws.Cells.Item(52, 13).Value = 19316
$ws.Cells.Item(52, 14).Value = '$/caja 20 kilos'
$ws.Cells.Item(52, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value = 966
$ws.Cells.Item(52, 17).Value = 20
$ws.Cells.Item(52, 18).Value = "Hortaliza"
